$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.251648
$ws.Range("H2").Value = 0.7549440000000001
$ws.Range("I2").Value = 0.6152542207125417
$ws.Range("J2").Value = 0.6152542207125418
$ws.Range("M2").Value = 0.032838
$ws.Range("N2").Value = 0.098514
$ws.Range("O2").Value = 0.007146324094219707
$ws.Range("P2").Value = 0.007146324094219707
$ws.Range("Q2").Value = 0.008263617024000001
$ws.Range("R2").Value = 0.074372553216
$ws.Range("S2").Value = 0.004396806061548406
$ws.Range("T2").Value = 0.004396806061548407

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.251648
$ws.Range("H3").Value = 0.7549440000000001
$ws.Range("I3").Value = 0.6152542207125417
$ws.Range("J3").Value = 0.6152542207125418
$ws.Range("M3").Value = 0.181585
$ws.Range("N3").Value = 0.544755
$ws.Range("O3").Value = 0.03951718316124263
$ws.Range("P3").Value = 0.03951718316124263
$ws.Range("Q3").Value = 0.04569550208000001
$ws.Range("R3").Value = 0.4112595187200001
$ws.Range("S3").Value = 0.02431311373062511
$ws.Range("T3").Value = 0.02431311373062511

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.251648
$ws.Range("H4").Value = 0.7549440000000001
$ws.Range("I4").Value = 0.6152542207125417
$ws.Range("J4").Value = 0.6152542207125418
$ws.Range("M4").Value = 3.814633
$ws.Range("N4").Value = 11.443899
$ws.Range("O4").Value = 0.8301542030119253
$ws.Range("P4").Value = 0.8301542030119253
$ws.Range("Q4").Value = 0.9599447651840002
$ws.Range("R4").Value = 8.639502886656
$ws.Range("S4").Value = 0.5107558772453432
$ws.Range("T4").Value = 0.5107558772453433

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.251648
$ws.Range("H5").Value = 0.7549440000000001
$ws.Range("I5").Value = 0.6152542207125417
$ws.Range("J5").Value = 0.6152542207125418
$ws.Range("M5").Value = 0.5660336666666667
$ws.Range("N5").Value = 1.698101
$ws.Range("O5").Value = 0.1231822897326124
$ws.Range("P5").Value = 0.1231822897326124
$ws.Range("Q5").Value = 0.1424412401493333
$ws.Range("R5").Value = 1.281971161344
$ws.Range("S5").Value = 0.07578842367502496
$ws.Range("T5").Value = 0.07578842367502497

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.017789
$ws.Range("H6").Value = 0.053367
$ws.Range("I6").Value = 0.04349232790348186
$ws.Range("J6").Value = 0.04349232790348186
$ws.Range("M6").Value = 0.032838
$ws.Range("N6").Value = 0.098514
$ws.Range("O6").Value = 0.007146324094219707
$ws.Range("P6").Value = 0.007146324094219707
$ws.Range("Q6").Value = 0.000584155182
$ws.Range("R6").Value = 0.005257396638
$ws.Range("S6").Value = 0.0003108102708103565
$ws.Range("T6").Value = 0.0003108102708103565

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.017789
$ws.Range("H7").Value = 0.053367
$ws.Range("I7").Value = 0.04349232790348186
$ws.Range("J7").Value = 0.04349232790348186
$ws.Range("M7").Value = 0.181585
$ws.Range("N7").Value = 0.544755
$ws.Range("O7").Value = 0.03951718316124263
$ws.Range("P7").Value = 0.03951718316124263
$ws.Range("Q7").Value = 0.003230215565
$ws.Range("R7").Value = 0.029071940085
$ws.Range("S7").Value = 0.001718694287870716
$ws.Range("T7").Value = 0.001718694287870716

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.017789
$ws.Range("H8").Value = 0.053367
$ws.Range("I8").Value = 0.04349232790348186
$ws.Range("J8").Value = 0.04349232790348186
$ws.Range("M8").Value = 3.814633
$ws.Range("N8").Value = 11.443899
$ws.Range("O8").Value = 0.8301542030119253
$ws.Range("P8").Value = 0.8301542030119253
$ws.Range("Q8").Value = 0.067858506437
$ws.Range("R8").Value = 0.610726557933
$ws.Range("S8").Value = 0.0361053388078483
$ws.Range("T8").Value = 0.03610533880784831

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.017789
$ws.Range("H9").Value = 0.053367
$ws.Range("I9").Value = 0.04349232790348186
$ws.Range("J9").Value = 0.04349232790348186
$ws.Range("M9").Value = 0.5660336666666667
$ws.Range("N9").Value = 1.698101
$ws.Range("O9").Value = 0.1231822897326124
$ws.Range("P9").Value = 0.1231822897326124
$ws.Range("Q9").Value = 0.01006917289633333
$ws.Range("R9").Value = 0.09062255606699998
$ws.Range("S9").Value = 0.005357484536952485
$ws.Range("T9").Value = 0.005357484536952485

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.1395776666666667
$ws.Range("H10").Value = 0.418733
$ws.Range("I10").Value = 0.3412534513839764
$ws.Range("J10").Value = 0.3412534513839764
$ws.Range("M10").Value = 0.032838
$ws.Range("N10").Value = 0.098514
$ws.Range("O10").Value = 0.007146324094219707
$ws.Range("P10").Value = 0.007146324094219707
$ws.Range("Q10").Value = 0.004583451418000001
$ws.Range("R10").Value = 0.041251062762
$ws.Range("S10").Value = 0.002438707761860944
$ws.Range("T10").Value = 0.002438707761860945

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.1395776666666667
$ws.Range("H11").Value = 0.418733
$ws.Range("I11").Value = 0.3412534513839764
$ws.Range("J11").Value = 0.3412534513839764
$ws.Range("M11").Value = 0.181585
$ws.Range("N11").Value = 0.544755
$ws.Range("O11").Value = 0.03951718316124263
$ws.Range("P11").Value = 0.03951718316124263
$ws.Range("Q11").Value = 0.02534521060166667
$ws.Range("R11").Value = 0.228106895415
$ws.Range("S11").Value = 0.0134853751427468
$ws.Range("T11").Value = 0.0134853751427468

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1395776666666667
$ws.Range("H12").Value = 0.418733
$ws.Range("I12").Value = 0.3412534513839764
$ws.Range("J12").Value = 0.3412534513839764
$ws.Range("M12").Value = 3.814633
$ws.Range("N12").Value = 11.443899
$ws.Range("O12").Value = 0.8301542030119253
$ws.Range("P12").Value = 0.8301542030119253
$ws.Range("Q12").Value = 0.5324375733296668
$ws.Range("R12").Value = 4.791938159967001
$ws.Range("S12").Value = 0.2832929869587338
$ws.Range("T12").Value = 0.2832929869587338

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1395776666666667
$ws.Range("H13").Value = 0.418733
$ws.Range("I13").Value = 0.3412534513839764
$ws.Range("J13").Value = 0.3412534513839764
$ws.Range("M13").Value = 0.5660336666666667
$ws.Range("N13").Value = 1.698101
$ws.Range("O13").Value = 0.1231822897326124
$ws.Range("P13").Value = 0.1231822897326124
$ws.Range("Q13").Value = 0.07900565844811112
$ws.Range("R13").Value = 0.711050926033
$ws.Range("S13").Value = 0.04203638152063495
$ws.Range("T13").Value = 0.04203638152063494
